$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 81: "check db in deploy:" -> "Reviewed Kelly's project. Mostly finished front end of Indie project."
$ws.Range("D81").Value = "Reviewed Kelly's project. Mostly finished front end of Indie project."

# Row 82: "Reviewed Kelly's project..." -> "Looked into why new stories weren't searchable on aws."
$ws.Range("D82").Value = "Looked into why new stories weren't searchable on aws."

# Row 83: "Looked into why..." -> revised text, hours 2.5 -> 7
$ws.Range("B83").Value = 7
$ws.Range("D83").Value = "Revised search jsp to hold inputs when search validation failed. Refactored code to reduce class and method size. Started testing methods to get chart data."
$ws.Rows("83").RowHeight = 30

# Row 84: new date/hours/task (used to be blank A84 only)
$ws.Range("A84").Value = 43595
$ws.Range("B84").Value = 2
$ws.Range("D84").Value = "Added data to database locally."

# Row 85: D85 ("10:30 - x") removed entirely, leave blank
$ws.Range("D85").Clear()

# Row 86: new text "Fri plan"
$ws.Range("D86").Value = "Fri plan"

# Row 87: new text
$ws.Range("D87").Value = "Make new users to demo a cluster of data; add via program."

# Row 88: A88 shared string ("check db in deploy:") removed -> blank; D88 new text
$ws.Range("A88").Value = ""
$ws.Range("D88").Value = "Make new users have interesting stories."

# Row 90: new text
$ws.Range("D90").Value = "Choose searches  to demo."

# Row 91: new text
$ws.Range("D91").Value = "1. retrieve cluster by career and family size"

# Row 93: new text
$ws.Range("D93").Value = "3. attempt retrieval with no results"

# Row 94: new text
$ws.Range("D94").Value = "a. orthodontist to demo failed api search"

# Row 95: new text
$ws.Range("D95").Value = "b. another one with no data cluster."

# Row 92: new text (added last, so it lands at the end of the shared string table)
$ws.Range("D92").Value = "2. retrieve whatever possible that requires 30% search"

$ws.Range("D93").Select()
